# Atualização de bases das ligas, do dia: 09-04-2024 às 22:40
#
# The underlying data rows got refreshed; in several places this resulted in
# two data rows trading places (everything but the leading running-counter
# in column A swapped between the two rows), plus a couple of odds values
# being refreshed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB, $firstCol, $lastCol) {
    # Swap the contents of columns $firstCol..$lastCol between two rows,
    # leaving column A (the running id) untouched.
    # Value2 (not Value) is used so numbers/strings round-trip as plain
    # variants instead of a wrapped Value object.
    $rangeA = $ws.Range($ws.Cells.Item($rowA, $firstCol), $ws.Cells.Item($rowA, $lastCol))
    $rangeB = $ws.Range($ws.Cells.Item($rowB, $firstCol), $ws.Cells.Item($rowB, $lastCol))

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# Columns B (2) through AC (29) swap between the row pairs below.
Swap-Rows $ws 38 39 2 29
Swap-Rows $ws 129 131 2 29
Swap-Rows $ws 192 193 2 29

# Row 231 just got two odds values refreshed.
$ws.Range("R231").Value2 = 1.85
$ws.Range("S231").Value2 = 1.95
